$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# Header E1: CPI -> logCPI100 (shared-string table reorders automatically,
# matching the removal of the old "CPI" entry and the append of "logCPI100")
$ws.Range("E1").Value = "logCPI100"

# Bulk-update column E (rows 2-301): raw CPI levels -> log(CPI)*100 series
$arr = New-Object 'object[,]' 300,1
$arr[0,0] = 214.45742076096161
$arr[1,0] = 214.67480136306398
$arr[2,0] = 214.76763242410985
$arr[3,0] = 214.89109931093563
$arr[4,0] = 215.01421618485588
$arr[5,0] = 215.16762308470479
$arr[6,0] = 215.32049000842846
$arr[7,0] = 215.50322287909705
$arr[8,0] = 215.65491513317815
$arr[9,0] = 215.74567681342259
$arr[10,0] = 215.89652603834102
$arr[11,0] = 216.04685311190374
$arr[12,0] = 216.1667412437736
$arr[13,0] = 216.2564406523019
$arr[14,0] = 216.40552918934515
$arr[15,0] = 216.52443261253106
$arr[16,0] = 216.67260555800519
$arr[17,0] = 216.76126727275303
$arr[18,0] = 216.90863574870227
$arr[19,0] = 216.99681739968923
$arr[20,0] = 217.05550585212086
$arr[21,0] = 217.26029312098598
$arr[22,0] = 217.37688231366499
$arr[23,0] = 217.49315935284423
$arr[24,0] = 217.63806922432701
$arr[25,0] = 217.78249718646819
$arr[26,0] = 217.84013415337552
$arr[27,0] = 218.0125875164054
$arr[28,0] = 218.12717715594613
$arr[29,0] = 218.21292140529982
$arr[30,0] = 218.26999033360428
$arr[31,0] = 218.41233542396711
$arr[32,0] = 218.44074854123204
$arr[33,0] = 218.5825359612962
$arr[34,0] = 218.72386198314788
$arr[35,0] = 218.83659260631484
$arr[36,0] = 218.89284837608534
$arr[37,0] = 218.9770956346874
$arr[38,0] = 219.11714557285586
$arr[39,0] = 219.17303933628563
$arr[40,0] = 219.28461151888419
$arr[41,0] = 219.42367487238292
$arr[42,0] = 219.50689964685898
$arr[43,0] = 219.61761850399733
$arr[44,0] = 219.72805581256191
$arr[45,0] = 219.81069988734015
$arr[46,0] = 219.92064791616576
$arr[47,0] = 219.9480914862356
$arr[48,0] = 220.13971243204514
$arr[49,0] = 220.24883170600936
$arr[50,0] = 220.38484637462346
$arr[51,0] = 220.52043639481445
$arr[52,0] = 220.60158767633445
$arr[53,0] = 220.70955404192182
$arr[54,0] = 220.79035303860513
$arr[55,0] = 220.89785172762535
$arr[56,0] = 221.00508498751373
$arr[57,0] = 221.13875529368588
$arr[58,0] = 221.21876044039576
$arr[59,0] = 221.24539610402758
$arr[60,0] = 221.40486794119414
$arr[61,0] = 221.4843848047698
$arr[62,0] = 221.59018132040319
$arr[63,0] = 221.64298308762511
$arr[64,0] = 221.74839442139063
$arr[65,0] = 221.85355052165278
$arr[66,0] = 221.93225084193369
$arr[67,0] = 222.01080880400551
$arr[68,0] = 222.14142378423384
$arr[69,0] = 222.219604630172
$arr[70,0] = 222.32362731029974
$arr[71,0] = 222.37554536572412
$arr[72,0] = 222.47919564926818
$arr[73,0] = 222.5567713439471
$arr[74,0] = 222.63420871636308
$arr[75,0] = 222.76296495710088
$arr[76,0] = 222.84003587030048
$arr[77,0] = 222.8913405994688
$arr[78,0] = 222.96818423176759
$arr[79,0] = 222.99376859079339
$arr[80,0] = 223.09595557485687
$arr[81,0] = 223.1979026831504
$arr[82,0] = 223.24878663529861
$arr[83,0] = 223.35037603411342
$arr[84,0] = 223.45172835126866
$arr[85,0] = 223.52758766870525
$arr[86,0] = 223.60331471176357
$arr[87,0] = 223.67890994092926
$arr[88,0] = 223.77949932739227
$arr[89,0] = 223.85478876813275
$arr[90,0] = 223.92994791268924
$arr[91,0] = 224.05492482826
$arr[92,0] = 224.10481506716442
$arr[93,0] = 224.15464805965485
$arr[94,0] = 224.25414282983843
$arr[95,0] = 224.4029589030022
$arr[96,0] = 224.45245115700837
$arr[97,0] = 224.45245115700837
$arr[98,0] = 224.4771761495295
$arr[99,0] = 224.62523122993221
$arr[100,0] = 224.67447097238414
$arr[101,0] = 224.69906992415497
$arr[102,0] = 224.82185611900749
$arr[103,0] = 224.87087356009175
$arr[104,0] = 224.99317566341949
$arr[105,0] = 225.06639194632436
$arr[106,0] = 225.13948500401045
$arr[107,0] = 225.2124552505644
$arr[108,0] = 225.35802895621831
$arr[109,0] = 225.38224387080734
$arr[110,0] = 225.52725051033059
$arr[111,0] = 225.59957267224019
$arr[112,0] = 225.6958152560932
$arr[113,0] = 225.79184503140587
$arr[114,0] = 225.88766293721312
$arr[115,0] = 225.98326990634834
$arr[116,0] = 226.07866686549761
$arr[117,0] = 226.150077319828
$arr[118,0] = 226.26883443016962
$arr[119,0] = 226.31624649622165
$arr[120,0] = 226.45817292380775
$arr[121,0] = 226.57609167176105
$arr[122,0] = 226.64668954402413
$arr[123,0] = 226.74064187529041
$arr[124,0] = 226.78754193188973
$arr[125,0] = 226.95129442179166
$arr[126,0] = 227.04459080179626
$arr[127,0] = 227.11443179490783
$arr[128,0] = 227.207378750001
$arr[129,0] = 227.27695865517595
$arr[130,0] = 227.43887955503789
$arr[131,0] = 227.50808984568587
$arr[132,0] = 227.57719001649312
$arr[133,0] = 227.66915288450397
$arr[134,0] = 227.6921132065774
$arr[135,0] = 227.80673308886628
$arr[136,0] = 227.87536009528287
$arr[137,0] = 227.92105126013951
$arr[138,0] = 227.98949800116381
$arr[139,0] = 228.1260687055013
$arr[140,0] = 228.17149700272958
$arr[141,0] = 228.21687783046417
$arr[142,0] = 228.30749747354716
$arr[143,0] = 228.35273648616936
$arr[144,0] = 228.4205067701794
$arr[145,0] = 228.44307338445194
$arr[146,0] = 228.44307338445194
$arr[147,0] = 228.44307338445194
$arr[148,0] = 228.53322276438846
$arr[149,0] = 228.55573090077738
$arr[150,0] = 228.64564697469828
$arr[151,0] = 228.69053529723749
$arr[152,0] = 228.71296207191111
$arr[153,0] = 228.7801729930226
$arr[154,0] = 228.7801729930226
$arr[155,0] = 228.82492255719859
$arr[156,0] = 228.91428359323331
$arr[157,0] = 228.98118391176214
$arr[158,0] = 229.11467617318854
$arr[159,0] = 229.20344359947364
$arr[160,0] = 229.26990030439299
$arr[161,0] = 229.35835134961167
$arr[162,0] = 229.40250940953226
$arr[163,0] = 229.42457161381182
$arr[164,0] = 229.55670999624789
$arr[165,0] = 229.64457942063964
$arr[166,0] = 229.73227142053028
$arr[167,0] = 229.79792441593622
$arr[168,0] = 229.88530764097069
$arr[169,0] = 229.97251539756371
$arr[170,0] = 230.12470886362112
$arr[171,0] = 230.14640731432996
$arr[172,0] = 230.21143769562008
$arr[173,0] = 230.23309286843991
$arr[174,0] = 230.29799367482494
$arr[175,0] = 230.34120705967419
$arr[176,0] = 230.38437748886545
$arr[177,0] = 230.53513694466238
$arr[178,0] = 230.64250275506873
$arr[179,0] = 230.70679506612984
$arr[180,0] = 230.79237036118818
$arr[181,0] = 230.87777736647212
$arr[182,0] = 231.02683666324478
$arr[183,0] = 231.13299523037932
$arr[184,0] = 231.26004392612595
$arr[185,0] = 231.36563466180314
$arr[186,0] = 231.44992279731517
$arr[187,0] = 231.55505344219048
$arr[188,0] = 231.63897510731951
$arr[189,0] = 231.72273491764201
$arr[190,0] = 231.76455432211586
$arr[191,0] = 231.8272080211627
$arr[192,0] = 231.93143040905122
$arr[193,0] = 232.04267206936368
$arr[194,0] = 232.10140076422809
$arr[195,0] = 232.1695757854668
$arr[196,0] = 232.23392261724848
$arr[197,0] = 232.3029222069585
$arr[198,0] = 232.38149770242572
$arr[199,0] = 232.45273201093482
$arr[200,0] = 232.5421241255612
$arr[201,0] = 232.64935713232035
$arr[202,0] = 232.76858095410242
$arr[203,0] = 232.87220105492631
$arr[204,0] = 232.99487888324455
$arr[205,0] = 233.0289961480791
$arr[206,0] = 233.12652916762261
$arr[207,0] = 233.15487605075239
$arr[208,0] = 233.23091623129523
$arr[209,0] = 233.32940856236922
$arr[210,0] = 233.43833736578478
$arr[211,0] = 233.52432078645523
$arr[212,0] = 233.58849641990184
$arr[213,0] = 233.60352387610965
$arr[214,0] = 233.63536489687849
$arr[215,0] = 233.6309606123844
$arr[216,0] = 233.71516519162606
$arr[217,0] = 233.80419230796969
$arr[218,0] = 233.89602220201789
$arr[219,0] = 233.98606976748005
$arr[220,0] = 234.02536974384475
$arr[221,0] = 234.06661630060336
$arr[222,0] = 234.09653519167199
$arr[223,0] = 234.1426610255788
$arr[224,0] = 234.2264726835412
$arr[225,0] = 234.34105633888782
$arr[226,0] = 234.37354226501347
$arr[227,0] = 234.41583598277788
$arr[228,0] = 234.36704702349712
$arr[229,0] = 234.38633308213625
$arr[230,0] = 234.39656302418265
$arr[231,0] = 234.4042338996656
$arr[232,0] = 234.43175921950595
$arr[233,0] = 234.47733423332375
$arr[234,0] = 234.51050319054701
$arr[235,0] = 234.53913764775129
$arr[236,0] = 234.5787240795091
$arr[237,0] = 234.60202792045575
$arr[238,0] = 234.6644362611959
$arr[239,0] = 234.70234611551905
$arr[240,0] = 234.79210342336216
$arr[241,0] = 234.87194844584357
$arr[242,0] = 234.91881332741991
$arr[243,0] = 234.97183991320463
$arr[244,0] = 235.05871784309807
$arr[245,0] = 235.15972740562304
$arr[246,0] = 235.26030974949239
$arr[247,0] = 235.38431694838047
$arr[248,0] = 235.43351354967089
$arr[249,0] = 235.50797106860628
$arr[250,0] = 235.58325818459807
$arr[251,0] = 235.68000093778244
$arr[252,0] = 235.7700493336767
$arr[253,0] = 235.79996053843115
$arr[254,0] = 235.88443884820859
$arr[255,0] = 235.96552791258028
$arr[256,0] = 236.02600867765625
$arr[257,0] = 236.10153866562862
$arr[258,0] = 236.16711852165596
$arr[259,0] = 236.21675725116239
$arr[260,0] = 236.29704048540233
$arr[261,0] = 236.36570990529458
$arr[262,0] = 236.42319578445461
$arr[263,0] = 236.49728908307401
$arr[264,0] = 236.59164520673266
$arr[265,0] = 236.65518254991588
$arr[266,0] = 236.69717826768567
$arr[267,0] = 236.70426687113516
$arr[268,0] = 236.73484652627604
$arr[269,0] = 236.80078052211746
$arr[270,0] = 236.8993085165246
$arr[271,0] = 236.98390124567572
$arr[272,0] = 237.05130895985923
$arr[273,0] = 237.09218408579514
$arr[274,0] = 237.17308101446525
$arr[275,0] = 237.24682807337803
$arr[276,0] = 237.28402280360217
$arr[277,0] = 237.32523123118398
$arr[278,0] = 237.40606269718532
$arr[279,0] = 237.48802635371354
$arr[280,0] = 237.57166403655506
$arr[281,0] = 237.62794181794632
$arr[282,0] = 237.69326408780759
$arr[283,0] = 237.73153594978905
$arr[284,0] = 237.80088605604232
$arr[285,0] = 237.87448344975135
$arr[286,0] = 237.92293508704518
$arr[287,0] = 237.94578114479447
$arr[288,0] = 237.98691000711477
$arr[289,0] = 238.05223746143253
$arr[290,0] = 238.15753153183201
$arr[291,0] = 238.26401053410655
$arr[292,0] = 238.3255087879065
$arr[293,0] = 238.39302055347207
$arr[294,0] = 238.48281361700873
$arr[295,0] = 238.52790891072496
$arr[296,0] = 238.61706682536231
$arr[297,0] = 238.69766941635987
$arr[298,0] = 238.78185694333865
$arr[299,0] = 238.83623394474779
$ws.Range("E2:E301").Value = $arr

# Update active selection to E2 (matches post-edit sheet view / drop frozen scroll)
[void]$ws.Range("E2").Select()

